$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-27 18:52:43"

$wsZhCn.Range("H3").Value = "2016-08-27 18:52:38"
$wsZhCn.Range("K3").Value = "2016-08-27 18:52:56"

$wsDeDe.Range("K3").Value = "2016-08-27 18:53:08"
